$wb = $excel.ActiveWorkbook

# --- Hoja1!A1: update the "Conversión del día" text with new exchange rates ---
$hoja1 = $wb.Worksheets.Item("Hoja1")
$hoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.27 = 12467.32 pesos`n✅ 12467.32 pesos = 3.25 = 969.03 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas sheet: update N10, O10, N12, O12 ---
$tasas = $wb.Worksheets.Item("tasas")
$tasas.Range("N10").Value = 306
$tasas.Range("O10").Value = 3815
$tasas.Range("N12").Value = 3834
$tasas.Range("O12").Value = 298
